$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.137.16'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.827.52'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.011'
$ws.Range("E4").Value = '  +0.82%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.88'
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("E6").Value = '  +0.64%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4707'
$ws.Range("E7").Value = '  +0.53%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07398'
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8801'
$ws.Range("E10").Value = '  +0.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.33'
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.934.30'
$ws.Range("E12").Value = '  +6.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07327'
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.29'
$ws.Range("E14").Value = '  +2.11%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.377'
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008698'
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("E19").Value = '  +0.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.705.20'
$ws.Range("E20").Value = '  +2.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.65'
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.235'
$ws.Range("E22").Value = '  -0.94%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.099.48'
$ws.Range("E24").Value = '  +2.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.879'
$ws.Range("E25").Value = '  -0.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.70'
$ws.Range("E26").Value = '  +0.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.44'
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.139'
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.181'
$ws.Range("E29").Value = '  -1.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.22'
$ws.Range("E30").Value = '  -0.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08938'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.167'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7417'
$ws.Range("E33").Value = '  -2.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.511'
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.934'
$ws.Range("E35").Value = '  +0.84%  '
$ws.Range("E36").Value = '  +0.79%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.088'
$ws.Range("E37").Value = '  -1.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05300'
$ws.Range("E38").Value = '  -0.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01947'
$ws.Range("E39").Value = '  +0.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.411'
$ws.Range("E40").Value = '  +1.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.933'
$ws.Range("E41").Value = '  -1.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.210'
$ws.Range("E42").Value = '  +0.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5252'
$ws.Range("E43").Value = '  -0.86%  '
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.387'
$ws.Range("E45").Value = '  -0.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4862'
$ws.Range("E46").Value = '  -0.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.39'
$ws.Range("E47").Value = '  -1.04%  '
$ws.Range("E48").Value = '  +0.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '104.18'
$ws.Range("E49").Value = '  +0.76%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.652'
$ws.Range("E50").Value = '  -0.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06301'
$ws.Range("E51").Value = '  +0.00%  '
